$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.998.54'
$ws.Range("E2").Value = '  +1.34%  '
$ws.Range("D3").Value = '2.052.60'
$ws.Range("E3").Value = '  -2.20%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.49%  '
$ws.Range("E6").Value = '  +2.84%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.62'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +15.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '60.65'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.379'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0780'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.14%  '
$ws.Range("E12").Value = '  +6.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.92'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.22%  '
$ws.Range("D14").Value = '2.350.99'
$ws.Range("E14").Value = '  -2.07%  '
$ws.Range("E15").Value = '  -1.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.82%  '
$ws.Range("D17").Value = '2.054.37'
$ws.Range("E17").Value = '  -1.95%  '
$ws.Range("D18").Value = '36.938.22'
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("D19").Value = '0.0₃0923'
$ws.Range("E19").Value = '  +11.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.66%  '
$ws.Range("E22").Value = '  +3.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.54%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  -1.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.76%  '
$ws.Range("E30").Value = '  +1.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.57'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0622'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.02%  '
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0872'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.70%  '
$ws.Range("E37").Value = '  -5.23%  '
$ws.Range("E38").Value = '  -5.25%  '
$ws.Range("E39").Value = '  +0.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.101'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +20.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +10.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0223'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("E43").Value = '  -2.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.23%  '
$ws.Range("E45").Value = '  +1.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +51.41%  '
$ws.Range("E47").Value = '  +6.56%  '
$ws.Range("B48").Value = 'Gas'
$ws.Range("C48").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '13.19'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -53.08%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.291.56'
$ws.Range("E49").Value = '  -2.63%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.92'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.75%  '
$ws.Range("E51").Value = '  +6.60%  '
